# Add new skill/event entry as row 34 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

$ws.Cells.Item($row, 1).Value2  = 43
$ws.Cells.Item($row, 2).Value2  = "06- Workshops, Bootcamps and Presentations"
$ws.Cells.Item($row, 3).Value2  = 45776
$ws.Cells.Item($row, 4).Value2  = "Agentic AI Innovation Session with AWS and Salesforce"
$ws.Cells.Item($row, 5).Value2  = "Agentic AI products are revolutionizing the tech landscape. Clients can use Agentic AI to deploy a digital workforce,  enhancing their workflows with efficiency and speed.  AWS and Salesforce explore  their Agentic AI products and perform a demo of a  use case Salesforce built for DOGE.  This is an ideation session for Booz Allen technical  and sales leader for use cases we can  co-develop and close with our partners. "
$ws.Cells.Item($row, 7).Value2  = "artificial intelligence"
$ws.Cells.Item($row, 8).Value2  = "cloud platform"
$ws.Cells.Item($row, 9).Value2  = "aws"
$ws.Cells.Item($row, 10).Value2 = "salesforce"
$ws.Cells.Item($row, 13).Value2 = "Booz Allen"

# Match the right-aligned numeric style used by the rest of column A (e_id)
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row uses the taller "first/last row" height, same as header row
$ws.Rows.Item($row).RowHeight = 15.75

# Leave the selection where the author left it
$ws.Range("E38").Select()
